$d = $word.ActiveDocument
$t = $d.Tables.Item($d.Tables.Count)

function Set-CellText($table, $row, $col, $text) {
    $cell = $table.Cell($row, $col)
    $rng = $cell.Range
    $rng.Text = $text
    $rng.Font.Italic = 1
}

# Row 4 = Part# 15 / Front Board Cover
Set-CellText $t 4 1  "15"
Set-CellText $t 4 2  "Front Board Cover"
Set-CellText $t 4 3  "1"
Set-CellText $t 4 4  ".11 thick, 6.38 length,`r13.5 width"
Set-CellText $t 4 5  "Mostly aesthetic and for platform for your feet"
Set-CellText $t 4 6  "Plastic"
Set-CellText $t 4 9  "Glossy"
Set-CellText $t 4 10 "Goes on top of the front part of the board"

# Row 5 = Part# 13 / Back Board Cover
Set-CellText $t 5 1  "13"
Set-CellText $t 5 2  "Back Board Cover"
Set-CellText $t 5 3  "1"
Set-CellText $t 5 4  ".1 thick, 6.25 length,`r15.2 width"
Set-CellText $t 5 5  "Mostly aesthetic and for platform for your feet"
Set-CellText $t 5 6  "Plastic"
Set-CellText $t 5 9  "Glossy"
Set-CellText $t 5 10 "Goes on top of the back part of the board for visual reasons"

# Move the _GoBack bookmark into the (empty) General Notes cell of row 5.
# A direct Range into a run-less (empty) paragraph cell cannot reliably host a
# new bookmark in this engine, so we temporarily place a placeholder character,
# anchor the bookmark right before it, then remove the placeholder again.
$notesCell = $t.Cell(5, 11)
$notesRng = $notesCell.Range
$notesRng.Text = "X"

$notesCell2 = $t.Cell(5, 11)
$notesRng2 = $notesCell2.Range
$anchor = $d.Range($notesRng2.Start, $notesRng2.Start)
$d.Bookmarks.Add("_GoBack", $anchor)

$placeholder = $d.Range($notesRng2.Start, $notesRng2.Start + 1)
$placeholder.Text = ""

Write-Host "Done"
